# A new observation was recorded for "Macroferia Regional de Talca - Pepino
# ensalada": insert a new data row right before the current row 112, pushing
# the existing rows 112-209 down to 113-210. Populate the newly inserted
# row 112 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(112).Insert()

$ws.Cells.Item(112, 1).Value = 5
$ws.Cells.Item(112, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value = "Maule"
$ws.Cells.Item(112, 4).Value = 44447
$ws.Cells.Item(112, 5).Value = 7
$ws.Cells.Item(112, 6).Value = 100112043
$ws.Cells.Item(112, 7).Value = "Pepino ensalada"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 300
$ws.Cells.Item(112, 11).Value = 16000
$ws.Cells.Item(112, 12).Value = 16000
$ws.Cells.Item(112, 13).Value = 16000
$ws.Cells.Item(112, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(112, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(112, 16).Value = 267
$ws.Cells.Item(112, 17).Value = 60
$ws.Cells.Item(112, 18).Value = "Hortaliza"
